$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert "TED" as the new row 71 -------------------------------------
# (pushes the existing row 71 "Yelp" and everything below it down by one)
$ws.Rows("71:71").Insert()
$ws.Range("A71").Value = "TED"
$ws.Range("B71").Value = "FF2B06"
$ws.Range("C71").Value = 9
$ws.Range("D71").Formula = "=MOD((C71+100),360)"
$ws.Range("E71").Value = 98
$ws.Range("F71").Value = 100

# --- Insert "Laravel" as the new row 73 ---------------------------------
# (row 72 is now "Yelp"; the new row goes right after it, pushing
# "Google+" and everything below it down by one more)
$ws.Rows("73:73").Insert()
$ws.Range("A73").Value = "Laravel"
$ws.Range("B73").Value = "FB502B"
$ws.Range("C73").Value = 7
$ws.Range("D73").Formula = "=MOD((C73+100),360)"
$ws.Range("E73").Value = 76
$ws.Range("F73").Value = 98

# --- View state (best effort) -------------------------------------------
$ws.Range("A98").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 56
